# Fixed workflow: the sensitivity sweep now starts at cutoff-step 5 instead
# of step 1, so the first four (Cutoff, Reaction_number) rows of each sweep
# are dropped and every surviving row is shifted up to take the data of the
# row four below it. Column A (the 0-based row index) is left untouched.
# Applies to both worksheets ("NBR" and "BAR").

$wb = $excel.ActiveWorkbook

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Pull columns B (Cutoff) and C (Reaction_number) from 4 rows further
    # down into rows 2..16, leaving column A's sequential index as-is.
    for ($r = 2; $r -le 16; $r++) {
        $srcRow = $r + 4
        $cutoff = $ws.Cells.Item($srcRow, 2).Value2
        $reactionNumber = $ws.Cells.Item($srcRow, 3).Value2
        $ws.Cells.Item($r, 2).Value = $cutoff
        $ws.Cells.Item($r, 3).Value = $reactionNumber
    }

    # The last 4 rows (17:20) are now duplicates of what we just copied
    # upward; remove them so the sheet - and its dimension - shrink to
    # A1:C16.
    $ws.Rows("17:20").Delete()
}
